# "BD actualizada y zip"
#
# On the "GD" sheet, column C ("Inicio de Operaciones" / the operation
# start-date column) is removed entirely. Deleting the whole column shifts
# the "Electricidad generada (MWh)" column (old D) into C, and the
# "Emisiones de GEI Reducidas" formula column (old E) into D - exactly as
# Excel does natively on a column delete (styles, merged cells, and
# formulas such as D7/1000*Factores!$D$12 -> C7/1000*Factores!$D$12 all
# shift automatically). The now-unreferenced "Inicio de Operaciones"
# shared string is dropped by Excel on save.
#
# The GD sheet also becomes the active sheet/tab, with F11 left selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GD")

$ws.Activate()
$ws.Columns("C").Delete()

$ws.Range("F11").Select() | Out-Null
